$d = $word.ActiveDocument

function Replace-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $d.Range($p.Range.Start, $p.Range.End)
    $payload = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($payload)
}

# --- 1. "The base Manifest.json file was created..." paragraph ---
$manifestInner = '<w:body><w:p><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/><w:t xml:space="preserve">The base </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>Manifest.json</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> file was created, it contains basic information such as the extension name, the version of the program, a description, browser action (the extension logo and what pops up when it is clicked), and finally the scripts that execute on webpages visited.</w:t></w:r></w:p></w:body>'

# --- 2. "The content.js file is the JavaScript..." paragraph ---
$contentJsInner = '<w:body><w:p><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:tab/><w:t>The content.js file is the JavaScript that is injected into webpages. I added a basic script that uses &#8220;.</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>getElementByTagName</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>&#8221; to get all link tags and store them as an HTML collection. I then console.log all the links to display them in the console.</w:t></w:r></w:p></w:body>'

# --- 3. "The Background File was created..." paragraph ---
$backgroundInner = '<w:body><w:p><w:r><w:tab/><w:t xml:space="preserve">The Background File was created, it contains code that should be completed in the background independent of the user or any webpages in the browser. It was initially created as a way to send links back and forth from the extension to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>VirusTotal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> servers. This is accomplished using JavaScript along with an HTML form and the post method.</w:t></w:r></w:p></w:body>'

# --- 4. "The PUT method..." paragraph + four new paragraphs (Weeks 11/15, 11/22, etc.) ---
$putInner = '<w:body>' + `
  '<w:p><w:r><w:tab/><w:t xml:space="preserve">The PUT method has successfully placed the strings I give it into the address bar; however, Virus Total does not recognize the URLs I give it. This week I will look deeper into what URLs </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>VirusTotal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> needs as well as the basic ability for my scripts to communicate with </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>eachother</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>' + `
  '<w:p><w:r><w:t>Week of 11/15/2020</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:tab/><w:t>I implemented a function to convert URLs to their ascii form for special character. I made sure to include the: and / symbols for URLs in addition to other special characters to include other unforeseen cases. The next step is to begin collecting data from virus total to be analyzed and changed.</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:t>Week of 11/22/2020</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:tab/><w:t xml:space="preserve">The previous function was removed and replaced with the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>encodeURIComponent</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> function. I also began working on analyzing the Virus Total results which lead me to find that Virus Total has implemented an API that allows for quick and easy communication and analysis. I intent to change the current code in order to utilize this API instead of the raw PUT method. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' + `
  '</w:body>'

# Apply edits from bottom to top so paragraph indices of not-yet-edited
# paragraphs stay valid (the last edit adds new paragraphs, shifting
# everything after it, but there is nothing after it here).
Replace-ParagraphXml 61 $putInner
Replace-ParagraphXml 52 $backgroundInner
Replace-ParagraphXml 49 $contentJsInner
Replace-ParagraphXml 47 $manifestInner

Write-Output "Edits applied"
